$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Time value for row 3 (Venue B) from "11am" to "1:30pm"
$ws.Range("C3").Value = "1:30pm"

# Update the selected/active cell to match the saved selection state
$ws.Range("C3").Select()
